$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated date (E3): 250418 -> 250420 ---
$ws.Range("E3").Value = 250420

# --- Row 8-10: "check from cad mm3" placeholders replaced with real computed volumes ---
$ws.Range("C8").Value  = "35931.574 mm3"
$ws.Range("C9").Value  = "9213.417 mm3"
$ws.Range("C10").Value = "9002.649 mm3"

# --- Row 16: gear motor spec note replaced with actual datasheet reference ---
$ws.Range("E16").Value = "GEARMOTOR 200 RPM 3-6V DC (Adafruit 3777??)"

# --- Row 19-20: jumper-cable quantities turned into real numbers ---
$ws.Range("D19").Value = 17
$ws.Range("D20").Value = 6

# --- Rows 21-23: notes expanded ---
$ws.Range("E21").Value = "link or something?? Markings on the bag??"
$ws.Range("E22").Value = "?? Markings on the bag??"
$ws.Range("E23").Value = "?? Which one??"

# --- New trailing notes under the table ---
$ws.Range("B47").Value = "(brackets for not yet implemented/ready components/materials)"
$ws.Range("B49").Value = "remember to add jumpers with power source"
$ws.Range("B50").Value = "'+ 3 male-female"
$ws.Range("B51").Value = "'+ 2 female-female"
$ws.Range("B52").Value = "maybe"

# --- Column E widened to fit the new, longer text ---
$ws.Columns("E").AutoFit() | Out-Null

# --- Selection moved from I30 to E3 (also resets the scrolled topLeftCell back to A1) ---
$ws.Range("E3").Select() | Out-Null
